$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for all existing data rows (2..398)
for ($i = 2; $i -le 398; $i++) {
    $ws.Cells.Item($i, 3).Value = 45204
}

# Row 398 gets an explicit row height written (ht="15" customHeight="1")
$ws.Rows.Item(398).RowHeight = 15

# Add the new row 399
$ws.Cells.Item(399, 1).Value = "A 47463-2023"
$ws.Cells.Item(399, 2).Value = 45202
$ws.Cells.Item(399, 3).Value = 45204
$ws.Cells.Item(399, 4).Value = "HALLANDS LÄN"
$ws.Cells.Item(399, 5).Value = "HALMSTAD"
$ws.Cells.Item(399, 6).Value = "Kyrkan"
$ws.Cells.Item(399, 7).Value = 7.7
$ws.Cells.Item(399, 8).Value = 0
$ws.Cells.Item(399, 9).Value = 0
$ws.Cells.Item(399, 10).Value = 0
$ws.Cells.Item(399, 11).Value = 0
$ws.Cells.Item(399, 12).Value = 0
$ws.Cells.Item(399, 13).Value = 0
$ws.Cells.Item(399, 14).Value = 0
$ws.Cells.Item(399, 15).Value = 0
$ws.Cells.Item(399, 16).Value = 0
$ws.Cells.Item(399, 17).Value = 0
$ws.Cells.Item(399, 18).WrapText = $true

# B and C on row 399 are dates, formatted like other rows (yyyy-mm-dd style)
$ws.Cells.Item(399, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(399, 3).NumberFormat = "YYYY-MM-DD"
